# Updated Suzie's Roth IRA and 401K dividend figures for December 2016 on the
# "Yearly" sheet. Dependent totals (row 15) and the cross-sheet rollups on the
# "All Time" sheet (rows 7 and 46) recalc automatically.

$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

$wsYearly.Range("E14").Value = 40.29
$wsYearly.Range("F14").Value = 49.22

# Update selections on both sheets to match the saved view state; select the
# "Yearly" sheet's cell first, then "All Time" last so "All Time" ends up the
# active/tab-selected sheet, matching the original workbook's selection.
$wsYearly.Range("I14").Select()
$wsAllTime.Range("K12").Select()
